{"js": "// Replace each two-digit multiplication equation in the document with its\n// updated value, per the commit diff. Every \"old\" equation string occurs\n// exactly once in the document, so a case-sensitive exact search+replace\n// for each pair is unambiguous and preserves all run/paragraph formatting.\nconst replacements = [\n  [\"65\u00d787=5655\", \"58\u00d741=2378\"],\n  [\"12\u00d777=924\", \"43\u00d754=2322\"],\n  [\"60\u00d758=3480\", \"22\u00d722=484\"],\n  [\"59\u00d796=5664\", \"45\u00d715=675\"],\n  [\"80\u00d718=1440\", \"68\u00d716=1088\"],\n  [\"42\u00d729=1218\", \"49\u00d794=4606\"],\n  [\"58\u00d715=870\", \"12\u00d733=396\"],\n  [\"22\u00d749=1078\", \"51\u00d766=3366\"],\n  [\"78\u00d758=4524\", \"22\u00d798=2156\"],\n  [\"77\u00d720=1540\", \"82\u00d750=4100\"],\n  [\"93\u00d781=7533\", \"12\u00d762=744\"],\n  [\"93\u00d716=1488\", \"85\u00d763=5355\"],\n  [\"56\u00d795=5320\", \"56\u00d717=952\"],\n  [\"98\u00d762=6076\", \"42\u00d770=2940\"],\n  [\"90\u00d788=7920\", \"62\u00d751=3162\"],\n  [\"91\u00d717=1547\", \"66\u00d721=1386\"],\n  [\"92\u00d776=6992\", \"39\u00d792=3588\"],\n  [\"84\u00d755=4620\", \"35\u00d794=3290\"],\n  [\"91\u00d725=2275\", \"52\u00d764=3328\"],\n  [\"64\u00d724=1536\", \"12\u00d744=528\"],\n  [\"15\u00d727=405\", \"86\u00d728=2408\"],\n  [\"42\u00d720=840\", \"62\u00d792=5704\"],\n  [\"53\u00d753=2809\", \"38\u00d737=1406\"],\n  [\"51\u00d717=867\", \"49\u00d781=3969\"],\n  [\"39\u00d741=1599\", \"98\u00d758=5684\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation with its updated value,\n# per the commit diff. Every \"old\" equation string occurs exactly once in\n# the document, so a case-sensitive whole-document Find/Replace for each\n# pair is unambiguous and preserves all run/paragraph formatting.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"65\u00d787=5655\", \"58\u00d741=2378\"),\n  @(\"12\u00d777=924\", \"43\u00d754=2322\"),\n  @(\"60\u00d758=3480\", \"22\u00d722=484\"),\n  @(\"59\u00d796=5664\", \"45\u00d715=675\"),\n  @(\"80\u00d718=1440\", \"68\u00d716=1088\"),\n  @(\"42\u00d729=1218\", \"49\u00d794=4606\"),\n  @(\"58\u00d715=870\", \"12\u00d733=396\"),\n  @(\"22\u00d749=1078\", \"51\u00d766=3366\"),\n  @(\"78\u00d758=4524\", \"22\u00d798=2156\"),\n  @(\"77\u00d720=1540\", \"82\u00d750=4100\"),\n  @(\"93\u00d781=7533\", \"12\u00d762=744\"),\n  @(\"93\u00d716=1488\", \"85\u00d763=5355\"),\n  @(\"56\u00d795=5320\", \"56\u00d717=952\"),\n  @(\"98\u00d762=6076\", \"42\u00d770=2940\"),\n  @(\"90\u00d788=7920\", \"62\u00d751=3162\"),\n  @(\"91\u00d717=1547\", \"66\u00d721=1386\"),\n  @(\"92\u00d776=6992\", \"39\u00d792=3588\"),\n  @(\"84\u00d755=4620\", \"35\u00d794=3290\"),\n  @(\"91\u00d725=2275\", \"52\u00d764=3328\"),\n  @(\"64\u00d724=1536\", \"12\u00d744=528\"),\n  @(\"15\u00d727=405\", \"86\u00d728=2408\"),\n  @(\"42\u00d720=840\", \"62\u00d792=5704\"),\n  @(\"53\u00d753=2809\", \"38\u00d737=1406\"),\n  @(\"51\u00d717=867\", \"49\u00d781=3969\"),\n  @(\"39\u00d741=1599\", \"98\u00d758=5684\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
